# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) values for the
# 416ae302-...-zh-cn.xlf / 416ae302-...-de-de.xlf handback rows
# (rows 2 and 5 on each locale sheet share the same timestamp text).

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-21 12:18:12"
$wsZhCn.Range("H2").Value = "2016-03-21 12:18:34"
$wsZhCn.Range("E5").Value = "2016-03-21 12:18:12"
$wsZhCn.Range("H5").Value = "2016-03-21 12:18:34"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-21 12:18:16"
$wsDeDe.Range("H2").Value = "2016-03-21 12:18:40"
$wsDeDe.Range("E5").Value = "2016-03-21 12:18:16"
$wsDeDe.Range("H5").Value = "2016-03-21 12:18:40"
